$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "through" date advanced from 2022-08-05 to 2022-08-06, so rename the
# sheet and update the August row label to match.
$ws.Name = "Through 2022-08-06"
$ws.Cells.Replace("August (through 08-05)", "August (through 08-06)", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Refreshed counts for the August row (row 9)
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 34
$ws.Range("H9").Value = 35
$ws.Range("I9").Value = 31

# Refreshed counts for the Total row (row 10)
$ws.Range("B10").Value = 169
$ws.Range("C10").Value = 309
$ws.Range("D10").Value = 482
$ws.Range("E10").Value = 438
$ws.Range("F10").Value = 313
$ws.Range("G10").Value = 655
$ws.Range("H10").Value = 945
$ws.Range("I10").Value = 1001
